$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.970.00"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.384.73"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.28"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.89"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.60"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.391"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").Value = "3.965.24"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.63"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "3.384.00"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "61.076.41"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.72"
$ws.Range("E19").Value = "  -4.98%  "
$ws.Range("E20").Value = "  -5.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.21"
$ws.Range("E21").Value = "  -4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.85"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("E25").Value = "  -4.31%  "
$ws.Range("D26").Value = "3.526.76"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.37"
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.43"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "3.416.47"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.27"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "2.455.34"
$ws.Range("E48").Value = "  -6.06%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.92"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.75"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("E51").Value = "  +2.26%  "
